$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")
$rng = $ws.Range("A51:G52")
$rng.Borders.Color = 0
$rng.Borders.LineStyle = 1
$ws.Range("A52").Value = "Linking_AutoUser"
